# Update column E values on Sheet1 (listado.xlsx) to mark the selected
# product row within each product block ("Creamos detalle del producto").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    1  = 1
    9  = 0
    12 = 1
    15 = 0
    17 = 1
    20 = 0
    21 = 1
    29 = 0
    30 = 1
    33 = 0
    34 = 1
    41 = 1
    42 = 1
    46 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Range("E$row").Value = $changes[$row]
}
